# Updating filtered feeds from workflow
# Adds a new "title" column (C) to the filtered-feeds sheet, populated with
# the article title corresponding to each row's link/keyword.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new column, styled like the other header cells (A1/B1).
$ws.Range("C1").Value = "title"
$ws.Range("A1").Copy()
$ws.Range("C1").PasteSpecial(-4122)

# Column width for C (~50.71 chars wide in the source file); the engine
# snaps ColumnWidth to 1/6-character increments, so use a value that lands
# on the closest achievable grid point.
$ws.Columns.Item(3).ColumnWidth = 49.86

# Title text for each data row (rows 2-13), aligned with the existing
# link (column A) / keyword (column B) values already in the sheet.
$titles = @(
  "FDA D.I.S.C.O. Burst Edition: FDA approvals of Augtyro (repotrectinib) for NTRK gene fusion-positive solid tumors and Krazati (adagrasib) for KRAS G12C-mutated colorectal cancer",
  "FDA D.I.S.C.O. Burst Edition: FDA approvals of Tecentriq (atezolizumab) for unresectable or metastatic alveolar soft part sarcoma, and Krazati (adagrasib) for KRAS G12C-mutated locally advanced or metastatic non-small cell lung cancer",
  "FDA D.I.S.C.O. Burst Edition: FDA approval of Vijoice (alpelisib) for adult and pediatric patients two years of age and older with severe manifestations of PIK3CA-related overgrowth spectrum who require systemic therapy",
  "FDA D.I.S.C.O. Burst Edition: FDA approvals of Verzenio (abemaciclib) for adjuvant treatment of adult patients with hormone receptor-positive, human epidermal growth factor receptor 2-negative, node-positive, early breast cancer, & Keytruda (pembrolizumab) for persistent, recurrent or metastatic cervical cancer whose tumors express PD-L1 (CPS ≥1)",
  "FDA D.I.S.C.O. Burst Edition: FDA approvals of Lumakras (sotorasib) for patients with KRAS G12C  mutated locally advanced or metastatic non-small cell lung cancer, and Truseltiq (infigratinib) for unresectable locally advanced or metastatic cholangiocarcinoma with a fibroblast growth factor receptor 2 fusion or other rearrangement",
  "FDA D.I.S.C.O. Burst Edition: FDA approvals of Trodelvy (sacituzumab govitecan) for locally advanced/metastatic urothelial cancer who received platinum-containing chemotherapy & either PD-1/PD-L1 inhibitor and Opdivo (nivolumab) in combination with chemotherapy for metastatic gastric cancer and esophageal adenocarcinoma",
  "FDA D.I.S.C.O. Burst Edition: Libtayo (cemiplimab-rwlc) for first-line treatment of patients with advanced NSCLC (locally advanced who are not candidates for surgical resection or definitive chemoradiation or metastatic) whose tumors have high PD-L1 expression with no eGFR, anaplastic lymphoma kinase or receptor tyrosine kinase aberrations",
  "Lumea to Incorporate Myriad Genetics Cancer Tests Into Digital Pathology Platform",
  "Lumea to Incorporate Myriad Genetics Cancer Tests Into Digital Pathology Platform",
  "Lumea to Incorporate Myriad Genetics Cancer Tests Into Digital Pathology Platform",
  "Natera MRD Study Results Demonstrate Potential to Ressurect Failed Adjuvant Drugs",
  "Lumea to Incorporate Myriad Genetics Cancer Tests Into Digital Pathology Platform"
)

for ($i = 0; $i -lt $titles.Length; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 3).Value = $titles[$i]
}
